# Insert a new weekly record at row 28 (shifting the existing rows 28..86
# down to 29..87), then populate the new row with the latest week's data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(28).Insert()

$ws.Cells.Item(28, 1).Value = 9
$ws.Cells.Item(28, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(28, 3).Value = "Metropolitana"
$ws.Cells.Item(28, 4).Value = 44623
$ws.Cells.Item(28, 5).Value = 13
$ws.Cells.Item(28, 6).Value = 100112005
$ws.Cells.Item(28, 7).Value = "Puerro"
$ws.Cells.Item(28, 8).Value = "Sin especificar"
$ws.Cells.Item(28, 9).Value = "Primera"
$ws.Cells.Item(28, 10).Value = 106
$ws.Cells.Item(28, 11).Value = 7000
$ws.Cells.Item(28, 12).Value = 8000
$ws.Cells.Item(28, 13).Value = 7500
$ws.Cells.Item(28, 14).Value = "$/paquete 20 unidades"
$ws.Cells.Item(28, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(28, 16).Value = 375
$ws.Cells.Item(28, 17).Value = 20
$ws.Cells.Item(28, 18).Value = "Hortaliza"
